$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the previously-empty R/S averages for the mod_exp (64 bit) benchmark
# blocks (rows 14, 20 and 26 hold AVERAGE-row cells that were blank).
$ws.Range("R14").Value = 15546.12012
$ws.Range("S14").Value = 10017.89941

$ws.Range("R20").Value = 15355.68945
$ws.Range("S20").Value = 9420.85059

$ws.Range("R26").Value = 20531.10156
$ws.Range("S26").Value = 12246.55957

# Scroll the view so column J is the left-most visible column, matching the
# author's saved window position.
$excel.ActiveWindow.ScrollColumn = 10

# Move the selection to Q9:Q26 (active cell Q9), replacing the previous
# Q27:S27 selection.
$ws.Range("Q9:Q26").Select() | Out-Null
